$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to stay text (source data stores them as
# literal strings, e.g. "11.28", not numbers) before writing values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.386.81"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.778.45"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "313.76"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.5327"
$ws.Range("E7").Value = "  +12.26%  "
$ws.Range("D8").Value = "0.3767"
$ws.Range("E8").Value = "  +9.18%  "
$ws.Range("D9").Value = "42.84"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "0.07413"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("D12").Value = "0.9995"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "20.71"
$ws.Range("E13").Value = "  +4.89%  "
$ws.Range("D14").Value = "6.102"
$ws.Range("E14").Value = "  +4.96%  "
$ws.Range("D15").Value = "1.782.00"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "6.991"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "89.67"
$ws.Range("E17").Value = "  +3.48%  "
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "0.06436"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D22").Value = "5.908"
$ws.Range("E22").Value = "  +5.51%  "
$ws.Range("D23").Value = "27.427.58"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("E24").Value = "  +4.49%  "
$ws.Range("D25").Value = "2.093"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "155.38"
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").Value = "20.23"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "2.371"
$ws.Range("E28").Value = "  +15.08%  "
$ws.Range("D29").Value = "1.988.91"
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("D30").Value = "121.27"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").Value = "1.082"
$ws.Range("E31").Value = "  +5.70%  "
$ws.Range("D32").Value = "0.1034"
$ws.Range("E32").Value = "  +12.84%  "
$ws.Range("D33").Value = "5.579"
$ws.Range("E33").Value = "  +5.30%  "
$ws.Range("D34").Value = "3.626"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "0.02260"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").Value = "0.05971"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "4.923"
$ws.Range("E37").Value = "  +4.89%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "11.28"
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("D39").Value = "0.2056"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "0.6131"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  +10.14%  "
$ws.Range("D42").Value = "1.429"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").Value = "1.132"
$ws.Range("E43").Value = "  +4.63%  "
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("D45").Value = "0.5789"
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("D46").Value = "3.624"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").Value = "121.35"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").Value = "1.897"
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").Value = "1.119"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "0.06728"
$ws.Range("D51").Value = "70.86"
$ws.Range("E51").Value = "  +2.58%  "
